$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.309.12'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '2.058.21'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.93'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.39%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.383'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.42'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('E11').Value = '  +0.99%  '
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').Value = '2.361.90'
$ws.Range('E13').Value = '  +1.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.773'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.17'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').Value = '2.065.93'
$ws.Range('E18').Value = '  +1.71%  '
$ws.Range('D19').Value = '37.258.09'
$ws.Range('E19').Value = '  +1.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.26'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +12.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.84%  '
$ws.Range('D22').Value = '0.0₃0812'
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '225.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.64%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.10%  '
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.20%  '
$ws.Range('E28').Value = '  +6.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.81'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.45%  '
$ws.Range('E30').Value = '  -4.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  -0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.50'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.21%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.34%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0616'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.54'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.05%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.29'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('B39').Value = 'THORChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.84'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.75'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.67'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +13.84%  '
$ws.Range('E42').Value = '  +0.81%  '
$ws.Range('D43').Value = '1.480.06'
$ws.Range('E43').Value = '  +1.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '96.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0924'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.18%  '
$ws.Range('E46').Value = '  +4.14%  '
$ws.Range('E47').Value = '  +2.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.54'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('E49').Value = '  +1.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.88%  '
$ws.Range('E51').Value = '  +1.70%  '
